# Optional rebate for LED — can be switched from database.json5
#
# 1) "${PB}" -> "${MPB}" in the Payback Period table cell (adds a
#    database-driven "M" placeholder letter in front of "PB").
# 2) Footnote 1: drop the stray grammar-checker markup
#    (<w:proofErr .../> pair around "National  Laboratory") so the
#    sentence reads as one continuous run of text.

$d = $word.ActiveDocument

# --- 1. ${PB} -> ${MPB} ------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("`${PB}", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    # $rng now spans the matched "${PB}" text; "P" sits two characters
    # in, right after the "${" — insert the new "M" immediately before it.
    $insPos = $rng.Start + 2
    $insRng = $d.Range($insPos, $insPos)
    $insRng.InsertBefore("M")
}

# --- 2. Clean up the footnote text -------------------------------------
$footnote = $d.Footnotes.Item(1)
$targetText = "Lumen Maintenance and Light Loss Factors: Consequences of " + `
    "Current Design Practices for LEDs, Pacific Northwest National  Laboratory"

# Re-assigning the exact same visible text is treated as a no-op by the
# host (it only rewrites the underlying runs/proofErr markup when the
# text actually changes), so bounce through a throwaway value first to
# force the footnote body to be rebuilt as a single clean run.
$footnote.Range.Text = "_tmp_"
$footnote.Range.Text = $targetText

Write-Output "done"
